$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -11
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -1
$ws.Range("F12").Value = -11
$ws.Range("F13").Value = -3
$ws.Range("F18").Value = -3
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = -5
